$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers I1 ("I0") and J1 ("IF"), copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-70: each entry is (I, J)
$data = @(
    @(7,7),
    @(7,7),
    @(5,5),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,9),
    @(8,8),
    @(6,6),
    @(8,8),
    @(7,7),
    @(6,6),
    @(7,7),
    @(6,6),
    @(6,6),
    @(8,8),
    @(6,6),
    @(5,5),
    @(8,8),
    @(6,6),
    @(6,6),
    @(6,7),
    @(7,7),
    @(7,8),
    @(5,5),
    @(6,7),
    @(6,6),
    @(8,9),
    @(7,7),
    @(6,6),
    @(8,8),
    @(5,6),
    @(7,7),
    @(7,7),
    @(5,6),
    @(9,9),
    @(9,9),
    @(6,7),
    @(8,8),
    @(8,8),
    @(6,8),
    @(9,9),
    @(6,7),
    @(7,7),
    @(7,8),
    @(5,6),
    @(7,7),
    @(8,9),
    @(8,8),
    @(7,8),
    @(6,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(8,9),
    @(10,10),
    @(7,8),
    @(6,7),
    @(5,6),
    @(1,1),
    @(7,7),
    @(7,7),
    @(6,6),
    @(6,6),
    @(7,7),
    @(5,5)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = 2 + $idx
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
